# Logged Week 15 and simulated Week 16
# Update the "H" (Home) row totals on both the OFF and DEF sheets to
# reflect the newly logged Week 15 data plus the simulated Week 16 data.

$wb = $excel.ActiveWorkbook

# --- OFF sheet ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 190
$wsOff.Range("C2").Value = 137
$wsOff.Range("D2").Value = 47
$wsOff.Range("E2").Value = 19

# --- DEF sheet ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 231
$wsDef.Range("C2").Value = 163
$wsDef.Range("D2").Value = 55
$wsDef.Range("G2").Value = 5
